# Update "想去人数" (wanna-go count) figures on the 展览 and 全部类型 sheets
# to match the freshly generated data output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 624
$wsExpo.Range("F6").Value = 14351
$wsExpo.Range("F7").Value = 16595
$wsExpo.Range("F8").Value = 16
$wsExpo.Range("F27").Value = 6747
$wsExpo.Range("F33").Value = 5757
$wsExpo.Range("F34").Value = 104
$wsExpo.Range("F37").Value = 4835

# Sheet "全部类型" (all types - aggregated)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 624
$wsAll.Range("F6").Value = 14351
$wsAll.Range("F7").Value = 16595
$wsAll.Range("F8").Value = 16
$wsAll.Range("F28").Value = 6747
$wsAll.Range("F36").Value = 5757
$wsAll.Range("F37").Value = 104
$wsAll.Range("F40").Value = 4835
